$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell that had its text tweaked first (keeps shared-string order
# consistent with the authored workbook).
$ws.Range("A6").Value = "Test +Inf and -Inf"

# Update the header row titles (order matches the order new shared strings
# were appended in the target workbook).
$ws.Range("D1").Value = "Assigned"
$ws.Range("B1").Value = "Time"
$ws.Range("A1").Value = "Task"
$ws.Range("E1").Value = "Actual Time"
$ws.Range("F1").Value = "Code Complete"
$ws.Range("G1").Value = "Reviewer"
$ws.Range("C1").Value = "Risk (1-5)"

# Resize columns A and F. The ColumnWidth COM property is offset by 5/6 of a
# character from the width stored in the OOXML <col> element, so compensate
# to land exactly on width="15" / width="14".
$ws.Columns("A").ColumnWidth = 14.166666666666666
$ws.Columns("F").ColumnWidth = 13.166666666666666

# Move the active selection to C2.
$ws.Range("C2").Select()
